# run: with adjusted supply side/IRA forecast
# Replace old CORRECTED CBO + SUPPLY SIDE folder with new one
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Federal Corporate Taxes Contribution (row 10, "historical" / current)
$ws.Range("I10").Value = 0.1067
$ws.Range("J10").Value = -0.0452
$ws.Range("K10").Value = 0.295
$ws.Range("L10").Value = 0.2573
$ws.Range("M10").Value = 0.3088
$ws.Range("N10").Value = 0.2796
$ws.Range("O10").Value = 0.1561

# Federal Corporate Taxes Contribution (row 11, "difference" / projection)
$ws.Range("P11").Value = -0.0145
$ws.Range("Q11").Value = -0.0203
$ws.Range("R11").Value = -0.0138
$ws.Range("S11").Value = -0.0083
$ws.Range("T11").Value = 0.0126
$ws.Range("U11").Value = 0.0046
$ws.Range("V11").Value = -0.005
$ws.Range("W11").Value = -0.0105
$ws.Range("X11").Value = -1.7941

# Fiscal Impact (row 30, "historical" / current)
$ws.Range("I30").Value = -2.2874
$ws.Range("J30").Value = -0.5607
$ws.Range("K30").Value = 0.0638
$ws.Range("L30").Value = -0.1884
$ws.Range("M30").Value = 0.5287
$ws.Range("N30").Value = 0.2312
$ws.Range("O30").Value = -0.2545

# Fiscal Impact (row 31, "difference" / projection)
$ws.Range("P31").Value = -0.6538
$ws.Range("Q31").Value = -0.2369
$ws.Range("R31").Value = -0.0086
$ws.Range("S31").Value = -0.5226
$ws.Range("T31").Value = -0.6951
$ws.Range("U31").Value = -0.4849
$ws.Range("V31").Value = -0.4247
$ws.Range("W31").Value = -0.0558
$ws.Range("X31").Value = -73.1408

# Federal Corporate Taxes Contribution (row 66, "historical" / difference vs previous)
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0

# Federal Corporate Taxes Contribution (row 67, "difference" / difference vs previous)
$ws.Range("P67").Value = 0.1058
$ws.Range("Q67").Value = 0.2178
$ws.Range("R67").Value = 0.2097
$ws.Range("S67").Value = 0.2755
$ws.Range("T67").Value = 0.2321
$ws.Range("U67").Value = 0.1752
$ws.Range("V67").Value = 0.0305
$ws.Range("W67").Value = 0.1036
$ws.Range("X67").Value = -1.3669

# Fiscal Impact (row 86, "historical" / difference vs previous)
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0

# Fiscal Impact (row 87, "difference" / difference vs previous)
$ws.Range("P87").Value = 0.02
$ws.Range("Q87").Value = 0.1261
$ws.Range("R87").Value = 0.1105
$ws.Range("S87").Value = 0.1636
$ws.Range("T87").Value = 0.1244
$ws.Range("U87").Value = 0.0793
$ws.Range("V87").Value = -0.0543
$ws.Range("W87").Value = 0.0352
$ws.Range("X87").Value = -2.4937
